# 20150603 - swap TIM1/TIM4 functions on AQ32 target
# (to better align with AQ32 pin usage.)
#
# Renames the "AQ32 Use" / "TauLabs Use" column headers to
# "AQ32 Definition" / "TauLabs Definition", and swaps which timer
# (TIM1 vs TIM4) channels 1-4 PWM outputs / PPM / RangeFinder are
# wired to, on the "Pins" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pins")

# Rename the three repeated header pairs (columns B/C, F/G, J/K).
$ws.Range("B1").Value = "AQ32 Definition"
$ws.Range("C1").Value = "TauLabs Definition"
$ws.Range("F1").Value = "AQ32 Definition"
$ws.Range("G1").Value = "TauLabs Definition"
$ws.Range("J1").Value = "AQ32 Definition"
$ws.Range("K1").Value = "TauLabs Definition"

# PE9/PE11/PE13/PE15 (J11/J13/J15/J16 = RX5..RX8 on TIM1 CH1..CH4) now also
# drive the PWM Out 9..12 TauLabs definitions that used to live on TIM4.
$ws.Range("K11").Value = "PWM Out 9 (TIM1 CH1)"
$ws.Range("K13").Value = "PWM Out 10 (TIM1 CH2)"
$ws.Range("K15").Value = "PWM Out 11 (TIM1 CH3)"
$ws.Range("K16").Value = "PWM Out 12 (TIM1 CH4)"

# PD12/PD13 (F31/F32 = RX1/RX2 on TIM4 CH1/CH2) now carry PPM and
# RangeFinder instead of PWM Out 9/10.
$ws.Range("G31").Value = "PPM (TIM4 CH1)"
$ws.Range("G32").Value = "RangeFinder (TIM4 CH2)"

# PD14/PD15 (F33/F34 = RX3/RX4 on TIM4 CH3/CH4) no longer have a TauLabs
# definition at all (their PWM Out 11/12 roles moved to TIM1 above).
$ws.Range("G33").ClearContents()
$ws.Range("G34").ClearContents()

# Leave the selection where the author left it when saving.
$ws.Range("K1").Select() | Out-Null
